$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# --- Update the mapping values (BUS v1 debit) ---
$ws.Range("E2").Value  = "C"
$ws.Range("B4").Value  = "C"
$ws.Range("D4").Value  = "P"
$ws.Range("F4").Value  = "E"
$ws.Range("H4").Value  = "C"
$ws.Range("F7").Value  = "E"
$ws.Range("C9").Value  = "PX"
$ws.Range("D9").Value  = "EX"
$ws.Range("E11").Value = "C"

# --- Center-align the data block (masse/cout durite) ---
$ws.Range("B2:H11").HorizontalAlignment = -4108

# --- Color-scale conditional formatting over the data block ---
$ws.Range("B2:H11").FormatConditions.AddColorScale(3)

# --- Activate the "map" sheet/tab and set its selection + zoom ---
$ws.Activate()
$ws.Range("E15").Select()
$excel.ActiveWindow.Zoom = 190
